$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts all existing rows 2..174 down to 3..175)
$ws.Rows(2).Insert()

# The inserted row picks up formatting copied from the row above (the bold
# header row) - clear it so the new data row matches the plain formatting
# used by every other data row.
$ws.Range("A2:D2").ClearFormats()

# Populate the newly inserted row with the new exposure site record
$ws.Range("A2").Value = "230 Rosanna Rd, Rosanna VIC 3084"
$ws.Range("B2").Value = -37.740508
$ws.Range("C2").Value = 145.075152
$ws.Range("D2").Value = "Banyule (C)"
